$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.088.13'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").Value = '1.821.16'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5915'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.89%  '
$ws.Range("E7").Value = '  +0.39%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2742'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.97%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06794'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.99'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07500'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.94%  '
$ws.Range("D12").Value = '1.833.24'
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.676'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6239'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000009417'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '74.59'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.43%  '
$ws.Range("D17").Value = '28.807.66'
$ws.Range("E17").Value = '  -1.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.427'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -9.33%  '
$ws.Range("E19").Value = '  +0.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '208.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -9.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.38'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.77%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.775'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.86%  '
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '154.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1269'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.786'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.27'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.84%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06482'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.411'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.434'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.713'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.680'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.678'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.051'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.530'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6319'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.750'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.442'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01706'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.37%  '
$ws.Range("D40").Value = '1.133.45'
$ws.Range("E40").Value = '  -8.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8697'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.005'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.41%  '
$ws.Range("D43").Value = '1.973.05'
$ws.Range("E43").Value = '  -0.65%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.64'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '60.22'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.19%  '
$ws.Range("E46").Value = '  -3.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.572'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.78%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05476'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4517'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.02%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.270'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.18%  '
$ws.Range("B51").Value = 'Frax'
$ws.Range("C51").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.011'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.86%  '
